$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "Price" column (D) holds values that look numeric (e.g.
# "487.59") but must stay plain text -- some prices use dotted thousands
# separators (e.g. "56.314.16") so the whole column is text-typed in the
# source data. Force text via NumberFormat "@" before assigning, then
# restore the default "Normal" style so no stray formatting is left
# behind (matches the original, un-styled cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- Simple price / volume(1h) refreshes -----------------------------
Set-TextValue $ws.Range("D2") '56.314.16'
$ws.Range("E2").Value = '  +4.95%  '
Set-TextValue $ws.Range("D3") '2.479.67'
$ws.Range("E3").Value = '  +2.55%  '
$ws.Range("E4").Value = '  -0.28%  '
Set-TextValue $ws.Range("D5") '487.59'
$ws.Range("E5").Value = '  +5.40%  '
Set-TextValue $ws.Range("D6") '146.58'
$ws.Range("E6").Value = '  +11.40%  '
$ws.Range("E7").Value = '  -0.21%  '
Set-TextValue $ws.Range("D8") '0.510'
$ws.Range("E8").Value = '  +5.01%  '
Set-TextValue $ws.Range("D9") '2.487.49'
$ws.Range("E9").Value = '  +2.61%  '
Set-TextValue $ws.Range("D10") '5.77'
$ws.Range("E10").Value = '  +8.42%  '
Set-TextValue $ws.Range("D11") '0.0966'
$ws.Range("E11").Value = '  +1.88%  '
Set-TextValue $ws.Range("D12") '0.330'
$ws.Range("E12").Value = '  +5.33%  '
$ws.Range("E13").Value = '  +1.49%  '
Set-TextValue $ws.Range("D14") '2.906.95'
$ws.Range("E14").Value = '  +2.23%  '
Set-TextValue $ws.Range("D15") '56.302.85'
$ws.Range("E15").Value = '  +4.76%  '
Set-TextValue $ws.Range("D16") '21.06'
$ws.Range("E16").Value = '  +7.18%  '
$ws.Range("E17").Value = '  +1.93%  '
Set-TextValue $ws.Range("D18") '2.482.43'
$ws.Range("E18").Value = '  +2.69%  '
Set-TextValue $ws.Range("D19") '4.51'
$ws.Range("E19").Value = '  +8.04%  '
Set-TextValue $ws.Range("D20") '10.05'
$ws.Range("E20").Value = '  +6.77%  '
Set-TextValue $ws.Range("D21") '317.26'
$ws.Range("E21").Value = '  +3.34%  '
Set-TextValue $ws.Range("D22") '1.00'
$ws.Range("E22").Value = '  -0.02%  '
Set-TextValue $ws.Range("D23") '5.79'
$ws.Range("E23").Value = '  +8.48%  '
Set-TextValue $ws.Range("D24") '58.39'
$ws.Range("E24").Value = '  +4.40%  '
Set-TextValue $ws.Range("D25") '0.410'
$ws.Range("E25").Value = '  +6.32%  '
$ws.Range("E26").Value = '  -1.05%  '
$ws.Range("E27").Value = '  +5.48%  '
Set-TextValue $ws.Range("D28") '2.586.35'
$ws.Range("E28").Value = '  +3.21%  '
Set-TextValue $ws.Range("D29") '7.61'
$ws.Range("E29").Value = '  +7.12%  '
Set-TextValue $ws.Range("D30") '0.0₃0791'
$ws.Range("E30").Value = '  +10.22%  '
$ws.Range("E31").Value = '  -0.18%  '
Set-TextValue $ws.Range("D32") '149.70'
$ws.Range("E32").Value = '  +2.44%  '
Set-TextValue $ws.Range("D33") '18.14'
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("E34").Value = '  +5.69%  '
Set-TextValue $ws.Range("D35") '5.19'
$ws.Range("E35").Value = '  +4.23%  '
$ws.Range("E36").Value = '  +7.97%  '
Set-TextValue $ws.Range("D37") '3.73'
$ws.Range("E37").Value = '  +5.86%  '
Set-TextValue $ws.Range("D38") '0.861'
$ws.Range("E38").Value = '  +7.38%  '
Set-TextValue $ws.Range("D39") '34.14'
$ws.Range("E39").Value = '  +4.14%  '
Set-TextValue $ws.Range("D40") '3.50'
$ws.Range("E40").Value = '  +7.68%  '
$ws.Range("E43").Value = '  +2.46%  '
Set-TextValue $ws.Range("D44") '1.33'
$ws.Range("E44").Value = '  +7.52%  '
Set-TextValue $ws.Range("D45") '4.77'
$ws.Range("E45").Value = '  +13.87%  '
Set-TextValue $ws.Range("D46") '0.0925'
$ws.Range("E46").Value = '  +6.32%  '
Set-TextValue $ws.Range("D47") '259.12'
$ws.Range("E47").Value = '  +15.46%  '
Set-TextValue $ws.Range("D50") '17.53'
$ws.Range("E50").Value = '  +6.50%  '
Set-TextValue $ws.Range("D51") '1.870.41'
$ws.Range("E51").Value = '  -3.29%  '

# --- Rank swap: Hedera now outranks FirstDigitalUSD (rows 41/42) -----
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range("D41") '0.0557'
$ws.Range("E41").Value = '  +6.22%  '
$ws.Range("B42").Value = 'FirstDigitalUSD'
$ws.Range("C42").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range("D42") '0.994'
$ws.Range("E42").Value = '  -0.15%  '

# --- Rank swap: VeChain now outranks WhiteBITCoin (rows 48/49) -------
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D48") '0.0228'
$ws.Range("E48").Value = '  +5.37%  '
$ws.Range("B49").Value = 'WhiteBITCoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range("D49") '10.17'
$ws.Range("E49").Value = '  +0.72%  '
